# [fix][asset] Fix asset bug to upload asset
#
# The "asset" upload template had a stale "机架" (rack) column and was
# missing a "设备类型" (device type) column. This also refreshes the
# sheet's tracked sort range and leaves the "asset" tab (not "part") as
# the active/selected sheet with A2 selected, matching the refreshed
# template.

$wb = $excel.ActiveWorkbook

$assetSheet = $wb.Worksheets.Item("asset")
$partSheet  = $wb.Worksheets.Item("part")

# Drop the obsolete "机架" column (column A); everything shifts one to the
# left (机柜 -> A, U位 -> B, 设备名称 -> C, 设备型号 -> D, ...).
$assetSheet.Columns.Item(1).Delete()

# Insert a new column in front of the old "设备型号" column (now D) and
# label it "设备类型", so the header row reads:
# 机柜, U位, 设备名称, 设备类型, 设备型号, 资产编号, ...
$assetSheet.Columns.Item(4).Insert()
$assetSheet.Range("D1").Value = "设备类型"

# The asset table grew by one column - refresh the sheet's remembered
# sort range/condition to match (now 28 columns wide: A..AB).
$sortObj = $assetSheet.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($assetSheet.Range("A1"))
$sortObj.SetRange($assetSheet.Range("A2:AB66"))
$sortObj.Apply()

# "part" no longer keeps its old B8 selection / active-tab status; "asset"
# becomes the active sheet with A2 selected instead.
$partSheet.Activate()
$partSheet.Range("A2").Select()

$assetSheet.Activate()
$assetSheet.Range("A2").Select()
